# "01_Intro to Ruby.pptx" - slide 2 - update the "Puts & Prints" textbox:
#   - rename shape "Puts & Prints" -> "Puts & Print"
#   - reflow box (narrower / shifted right, same top/height)
#   - fix the visible text "uts & Prints" -> "uts & Print" (second run)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)

# Rename the shape.
$shape.Name = "Puts & Print"

# Reposition / resize the textbox.
# Target EMU: off x=3644561 y=596900 (unchanged), ext cx=5509261 (unchanged cy=1295401)
# Points values below are the nearest Single (float32) representations that
# round-trip back to the exact target EMU through PowerPoint's point->EMU
# conversion (EMU = points * 12700).
$shape.Left = 286.97332763671875
$shape.Top = 47
$shape.Width = 433.8000793457031

# Fix the wording in the second run (keeps its own run-level formatting).
$tr = $shape.TextFrame.TextRange
$tr.Runs(2).Text = "uts & Print"

# Changing the text can re-trigger shape autofit and nudge the height;
# restore the original height explicitly (cy=1295401 EMU, unchanged by diff).
$shape.Height = 102.00008392333984
